# Update the "Förändrad" (changed) date column (C) for all existing data
# rows (2-125) from 2023-09-13 (45182) to 2023-09-15 (45184), and append a
# new data row (126) for case "A 43335-2023".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C holds the "Förändrad" date for every existing record (rows 2..125).
$ws.Range("C2:C125").Value = 45184

# Give row 125 an explicit row height, matching the rest of the sheet
# (all other rows already carry ht="15" customHeight="1").
$ws.Rows.Item(125).RowHeight = 15

# Append the new record as row 126.
$row = 126
$ws.Cells.Item($row, 1).Value = "A 43335-2023"
$ws.Cells.Item($row, 2).Value = 45183
$ws.Cells.Item($row, 3).Value = 45184
$ws.Cells.Item($row, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item($row, 5).Value = "SJÖBO"
$ws.Cells.Item($row, 6).Value = "Kommuner"
$ws.Cells.Item($row, 7).Value = 4.9
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 0
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0

# Match the existing date formatting (style index 1, numFmt "YYYY-MM-DD")
# used by columns B and C throughout the table.
$ws.Range("B126:C126").NumberFormat = "YYYY-MM-DD"

# Match the wrap-text style used for column R ("Artnamn") throughout.
$ws.Range("R126").WrapText = $true
